# Weekly update: insert a new daily price record for Ciboulette
# (Vega Central Mapocho de Santiago) right before the existing row 275,
# pushing every following record down by one row.
#
# The new row mirrors the fixed/contextual columns of the row it is
# inserted in front of (mercado, region, codigos, calidad, unidad,
# origen, clasificacion, ...) while carrying fresh values for the date
# (D), volumen (J), precio minimo/maximo/promedio (K/L/M) and precio
# $/Kg (P).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$insertAt = 275

# Shift rows 275:363 down to 276:364, leaving a blank row 275 behind.
$ws.Rows.Item($insertAt).Insert()

# The record that used to live at $insertAt now sits one row below;
# copy its contextual columns into the freshly inserted row.
$sourceRow = $insertAt + 1
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item($insertAt, $col).Value2 = $ws.Cells.Item($sourceRow, $col).Value2
}

# Overwrite the columns that actually carry new data for this record.
$ws.Cells.Item($insertAt, 4).Value2 = 44627    # D - Fecha
$ws.Cells.Item($insertAt, 10).Value2 = 106     # J - Volumen
$ws.Cells.Item($insertAt, 11).Value2 = 1500    # K - Precio minimo
$ws.Cells.Item($insertAt, 12).Value2 = 1800    # L - Precio maximo
$ws.Cells.Item($insertAt, 13).Value2 = 1650    # M - Precio promedio ponderado
$ws.Cells.Item($insertAt, 16).Value2 = 550     # P - Precio $/Kg
